# Applies the "Updated cryptos list" refresh (coin prices / 1h volume %,
# plus a handful of row-pair reorderings) described by the commit diff.
#
# D (Price) and E (Volume(1h)) columns hold text-formatted numbers/percentages
# (e.g. "60.716.27", "  -0.26%  ") in the source workbook (t="inlineStr").
# Assigning a bare numeric-looking string via .Value lets the COM layer
# coerce it to a real number (dropping formatting like trailing zeros, the
# thousands dots, or the "%" text). Prefixing with a single quote forces
# Excel to keep it as literal text, matching the original cell type; the
# follow-up `.Style = "Normal"` strips the quote-prefix formatting Excel
# applies so the cell's style stays identical to before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'60.716.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.26%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'2.641.27"
$ws.Range("D3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'578.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.23%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'143.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.15%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.18%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.599"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.47%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("E9").Value = "'  +0.65%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'0.106"
$ws.Range("D10").Style = "Normal"

# Row 11
$ws.Range("E11").Value = "'  +1.83%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("E12").Value = "'  -1.21%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'3.110.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.99%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'26.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +11.58%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'60.694.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.26%  "
$ws.Range("E15").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'2.653.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.07%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'11.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.07%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("E19").Value = "'  +0.90%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'349.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.29%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("E21").Value = "'  -1.53%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("E22").Value = "'  +0.06%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("E23").Value = "'  +1.69%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'63.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.94%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").Value = "'0.161"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.51%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'0.989"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.75%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("E27").Value = "'  +3.45%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("E28").Value = "'  +9.51%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("E29").Value = "'  +0.35%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("E30").Value = "'  +6.85%  "
$ws.Range("E30").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'163.43"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.22%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'19.87"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.42%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("E34").Value = "'  +7.18%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("E35").Value = "'  +3.57%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("E36").Value = "'  +7.02%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").Value = "'1.66"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.97%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "'339.15"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +10.04%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("E39").Value = "'  +3.98%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("E40").Value = "'  +6.55%  "
$ws.Range("E40").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'5.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.78%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("E43").Value = "'  +2.14%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'20.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.59%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("E45").Value = "'  +2.14%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'132.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.94%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0248"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.11%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'0.0995"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.00%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = "'  +0.46%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.39%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'2.086.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.83%  "
$ws.Range("E51").Style = "Normal"
